$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: was an empty styled cell, now holds shared string "test@gmail.com"
$ws.Range("G2").Style = "Normal"
$ws.Range("G2").Value = "test@gmail.com"

# H2: was "No", now "Sí"
$ws.Range("H2").Value = "Sí"

# Update the active selection to G9
$ws.Range("G9").Select()
